$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B120").Value = "SingleUseId258"
$ws.Range("C120").Value = "Large"
$ws.Range("D120").Value = "Left"
$ws.Range("E120").Value = "LTR"
$ws.Range("F120").Value = "<value>"

$ws.Range("B121").Value = "SingleUseId264"
$ws.Range("C121").Value = "Large"
$ws.Range("D121").Value = "Left"
$ws.Range("E121").Value = "LTR"
$ws.Range("F121").Value = "Info"

$ws.Range("B122").Value = "SingleUseId265"
$ws.Range("C122").Value = "Default"
$ws.Range("D122").Value = "Left"
$ws.Range("E122").Value = "LTR"
$ws.Range("F122").Value = "CPU Usage: <value>%"

$ws.Range("B123").Value = "SingleUseId266"
$ws.Range("C123").Value = "Default"
$ws.Range("D123").Value = "Left"
$ws.Range("E123").Value = "LTR"
$ws.Range("F123").Value = "0"

$ws.Range("B124").Value = "SingleUseId267"
$ws.Range("C124").Value = "Default"
$ws.Range("D124").Value = "Center"
$ws.Range("E124").Value = "LTR"
$ws.Range("F124").Value = "Back"

$ws.Range("B125").Value = "SingleUseId268"
$ws.Range("C125").Value = "Small"
$ws.Range("D125").Value = "Left"
$ws.Range("E125").Value = "LTR"
$ws.Range("F125").Value = "<value>"

$ws.Range("B126").Value = "SingleUseId269"
$ws.Range("C126").Value = "Small"
$ws.Range("D126").Value = "Left"
$ws.Range("E126").Value = "LTR"
$ws.Range("F126").Value = "<value>"

$ws.Range("B127").Value = "SingleUseId270"
$ws.Range("C127").Value = "Small"
$ws.Range("D127").Value = "Left"
$ws.Range("E127").Value = "LTR"
$ws.Range("F127").Value = "1000"

$ws.Range("B128").Value = "SingleUseId271"
$ws.Range("C128").Value = "Small"
$ws.Range("D128").Value = "Left"
$ws.Range("E128").Value = "LTR"
$ws.Range("F128").Value = "-1000"

$ws.Range("B129").Value = "SingleUseId272"
$ws.Range("C129").Value = "Small"
$ws.Range("D129").Value = "Left"
$ws.Range("E129").Value = "LTR"
$ws.Range("F129").Value = "<value>"

$ws.Range("B130").Value = "SingleUseId273"
$ws.Range("C130").Value = "Small"
$ws.Range("D130").Value = "Left"
$ws.Range("E130").Value = "LTR"
$ws.Range("F130").Value = "0"

$ws.Range("B131").Value = "SingleUseId274"
$ws.Range("C131").Value = "Small"
$ws.Range("D131").Value = "Left"
$ws.Range("E131").Value = "LTR"
$ws.Range("F131").Value = "<value>"

$ws.Range("B132").Value = "SingleUseId275"
$ws.Range("C132").Value = "Small"
$ws.Range("D132").Value = "Left"
$ws.Range("E132").Value = "LTR"
$ws.Range("F132").Value = "720"

$ws.Range("B133").Value = "SingleUseId276"
$ws.Range("C133").Value = "Default"
$ws.Range("D133").Value = "Center"
$ws.Range("E133").Value = "LTR"
$ws.Range("F133").Value = "Signals"

$ws.Range("B134").Value = "SingleUseId285"
$ws.Range("C134").Value = "Large"
$ws.Range("D134").Value = "Left"
$ws.Range("E134").Value = "LTR"
$ws.Range("F134").Value = "Signals"

$ws.Range("B135").Value = "SingleUseId286"
$ws.Range("C135").Value = "Default"
$ws.Range("D135").Value = "Left"
$ws.Range("E135").Value = "LTR"
$ws.Range("F135").Value = "CPU Usage: <value>%"

$ws.Range("B136").Value = "SingleUseId287"
$ws.Range("C136").Value = "Default"
$ws.Range("D136").Value = "Left"
$ws.Range("E136").Value = "LTR"
$ws.Range("F136").Value = "0"

$ws.Range("B137").Value = "SingleUseId288"
$ws.Range("C137").Value = "Default"
$ws.Range("D137").Value = "Center"
$ws.Range("E137").Value = "LTR"
$ws.Range("F137").Value = "Back"

$ws.Range("B138").Value = "SingleUseId289"
$ws.Range("C138").Value = "Large"
$ws.Range("D138").Value = "Left"
$ws.Range("E138").Value = "LTR"
$ws.Range("F138").Value = "Time range:<value>"

$ws.Range("B139").Value = "SingleUseId290"
$ws.Range("C139").Value = "Large"
$ws.Range("D139").Value = "Left"
$ws.Range("E139").Value = "LTR"
$ws.Range("F139").Value = "Y-axis min:<value>"

$ws.Range("B140").Value = "SingleUseId292"
$ws.Range("C140").Value = "Large"
$ws.Range("D140").Value = "Left"
$ws.Range("E140").Value = "LTR"
$ws.Range("F140").Value = "Y-axis max:<value>"

$ws.Range("B141").Value = "SingleUseId293"
$ws.Range("C141").Value = "Large"
$ws.Range("D141").Value = "Left"
$ws.Range("E141").Value = "LTR"
$ws.Range("F141").Value = "-1"

$ws.Range("B142").Value = "SingleUseId294"
$ws.Range("C142").Value = "Large"
$ws.Range("D142").Value = "Left"
$ws.Range("E142").Value = "LTR"
$ws.Range("F142").Value = "1"

$ws.Range("B143").Value = "SingleUseId295"
$ws.Range("C143").Value = "Large"
$ws.Range("D143").Value = "Left"
$ws.Range("E143").Value = "LTR"
$ws.Range("F143").Value = "720"

$ws.Range("B144").Value = "SingleUseId296"
$ws.Range("C144").Value = "Default"
$ws.Range("D144").Value = "Left"
$ws.Range("E144").Value = "LTR"
$ws.Range("F144").Value = "Auto Y range"

$ws.Range("B145").Value = "SingleUseId305"
$ws.Range("C145").Value = "Large"
$ws.Range("D145").Value = "Left"
$ws.Range("E145").Value = "LTR"
$ws.Range("F145").Value = "None"

$ws.Range("B146").Value = "SingleUseId306"
$ws.Range("C146").Value = "Large"
$ws.Range("D146").Value = "Left"
$ws.Range("E146").Value = "LTR"
$ws.Range("F146").Value = "<value>"

$ws.Range("B147").Value = "SingleUseId307"
$ws.Range("C147").Value = "Large"
$ws.Range("D147").Value = "Left"
$ws.Range("E147").Value = "LTR"
$ws.Range("F147").Value = "None"

$ws.Range("B148").Value = "SingleUseId308"
$ws.Range("C148").Value = "Large"
$ws.Range("D148").Value = "Left"
$ws.Range("E148").Value = "LTR"
$ws.Range("F148").Value = "<value>"

$ws.Range("B149").Value = "SingleUseId309"
$ws.Range("C149").Value = "Large"
$ws.Range("D149").Value = "Left"
$ws.Range("E149").Value = "LTR"
$ws.Range("F149").Value = "None"

$ws.Range("B150").Value = "SingleUseId310"
$ws.Range("C150").Value = "Large"
$ws.Range("D150").Value = "Left"
$ws.Range("E150").Value = "LTR"
$ws.Range("F150").Value = "<value>"

$ws.Range("B151").Value = "SingleUseId311"
$ws.Range("C151").Value = "Large"
$ws.Range("D151").Value = "Left"
$ws.Range("E151").Value = "LTR"
$ws.Range("F151").Value = "None"

$ws.Range("B152").Value = "SingleUseId312"
$ws.Range("C152").Value = "Large"
$ws.Range("D152").Value = "Left"
$ws.Range("E152").Value = "LTR"
$ws.Range("F152").Value = "<value>"

$ws.Range("B153").Value = "SingleUseId313"
$ws.Range("C153").Value = "Large"
$ws.Range("D153").Value = "Left"
$ws.Range("E153").Value = "LTR"
$ws.Range("F153").Value = "None"

$ws.Range("B154").Value = "SingleUseId314"
$ws.Range("C154").Value = "Large"
$ws.Range("D154").Value = "Left"
$ws.Range("E154").Value = "LTR"
$ws.Range("F154").Value = "<value>"

$ws.Range("B155").Value = "SingleUseId315"
$ws.Range("C155").Value = "Large"
$ws.Range("D155").Value = "Left"
$ws.Range("E155").Value = "LTR"
$ws.Range("F155").Value = "None"

$ws.Range("B156").Value = "SingleUseId316"
$ws.Range("C156").Value = "Large"
$ws.Range("D156").Value = "Left"
$ws.Range("E156").Value = "LTR"
$ws.Range("F156").Value = "<value>"

$ws.Range("B157").Value = "SingleUseId317"
$ws.Range("C157").Value = "Large"
$ws.Range("D157").Value = "Left"
$ws.Range("E157").Value = "LTR"
$ws.Range("F157").Value = "None"

$ws.Range("B158").Value = "SingleUseId318"
$ws.Range("C158").Value = "Large"
$ws.Range("D158").Value = "Left"
$ws.Range("E158").Value = "LTR"
$ws.Range("F158").Value = "<value>"

$ws.Range("B159").Value = "SingleUseId319"
$ws.Range("C159").Value = "Large"
$ws.Range("D159").Value = "Left"
$ws.Range("E159").Value = "LTR"
$ws.Range("F159").Value = "None"

$ws.Range("B160").Value = "SingleUseId320"
$ws.Range("C160").Value = "Large"
$ws.Range("D160").Value = "Left"
$ws.Range("E160").Value = "LTR"
$ws.Range("F160").Value = "<value>"

$ws.Range("B161").Value = "SingleUseId321"
$ws.Range("C161").Value = "Large"
$ws.Range("D161").Value = "Left"
$ws.Range("E161").Value = "LTR"
$ws.Range("F161").Value = "None"

$ws.Range("B162").Value = "SingleUseId322"
$ws.Range("C162").Value = "Large"
$ws.Range("D162").Value = "Left"
$ws.Range("E162").Value = "LTR"
$ws.Range("F162").Value = "<value>"

$ws.Range("B163").Value = "SingleUseId323"
$ws.Range("C163").Value = "Large"
$ws.Range("D163").Value = "Left"
$ws.Range("E163").Value = "LTR"
$ws.Range("F163").Value = "None"

$ws.Range("B164").Value = "SingleUseId324"
$ws.Range("C164").Value = "Large"
$ws.Range("D164").Value = "Left"
$ws.Range("E164").Value = "LTR"
$ws.Range("F164").Value = "<value>"

$ws.Range("B165").Value = "SingleUseId325"
$ws.Range("C165").Value = "Large"
$ws.Range("D165").Value = "Left"
$ws.Range("E165").Value = "LTR"
$ws.Range("F165").Value = "None"

$ws.Range("B166").Value = "SingleUseId326"
$ws.Range("C166").Value = "Large"
$ws.Range("D166").Value = "Left"
$ws.Range("E166").Value = "LTR"
$ws.Range("F166").Value = "<value>"

$ws.Range("B167").Value = "SingleUseId327"
$ws.Range("C167").Value = "Large"
$ws.Range("D167").Value = "Left"
$ws.Range("E167").Value = "LTR"
$ws.Range("F167").Value = "None"

$ws.Range("B168").Value = "SingleUseId337"
$ws.Range("C168").Value = "Large"
$ws.Range("D168").Value = "Left"
$ws.Range("E168").Value = "LTR"
$ws.Range("F168").Value = "None"

$ws.Range("B169").Value = "SingleUseId338"
$ws.Range("C169").Value = "Large"
$ws.Range("D169").Value = "Left"
$ws.Range("E169").Value = "LTR"
$ws.Range("F169").Value = "<value>"

$ws.Range("B170").Value = "SingleUseId339"
$ws.Range("C170").Value = "Large"
$ws.Range("D170").Value = "Left"
$ws.Range("E170").Value = "LTR"
$ws.Range("F170").Value = "None"

$ws.Range("B171").Value = "SingleUseId340"
$ws.Range("C171").Value = "Large"
$ws.Range("D171").Value = "Left"
$ws.Range("E171").Value = "LTR"
$ws.Range("F171").Value = "<value>"

$ws.Range("B172").Value = "SingleUseId341"
$ws.Range("C172").Value = "Large"
$ws.Range("D172").Value = "Left"
$ws.Range("E172").Value = "LTR"
$ws.Range("F172").Value = "None"

$ws.Range("B173").Value = "SingleUseId342"
$ws.Range("C173").Value = "Large"
$ws.Range("D173").Value = "Left"
$ws.Range("E173").Value = "LTR"
$ws.Range("F173").Value = "<value>"

$ws.Range("B174").Value = "SingleUseId343"
$ws.Range("C174").Value = "Large"
$ws.Range("D174").Value = "Left"
$ws.Range("E174").Value = "LTR"
$ws.Range("F174").Value = "None"

$ws.Range("B175").Value = "SingleUseId344"
$ws.Range("C175").Value = "Large"
$ws.Range("D175").Value = "Left"
$ws.Range("E175").Value = "LTR"
$ws.Range("F175").Value = "<value>"

$ws.Range("B176").Value = "SingleUseId345"
$ws.Range("C176").Value = "Large"
$ws.Range("D176").Value = "Left"
$ws.Range("E176").Value = "LTR"
$ws.Range("F176").Value = "None"

$ws.Range("B177").Value = "SingleUseId346"
$ws.Range("C177").Value = "Large"
$ws.Range("D177").Value = "Left"
$ws.Range("E177").Value = "LTR"
$ws.Range("F177").Value = "<value>"

$ws.Range("B178").Value = "SingleUseId347"
$ws.Range("C178").Value = "Large"
$ws.Range("D178").Value = "Left"
$ws.Range("E178").Value = "LTR"
$ws.Range("F178").Value = "None"

$ws.Range("B179").Value = "SingleUseId348"
$ws.Range("C179").Value = "Large"
$ws.Range("D179").Value = "Left"
$ws.Range("E179").Value = "LTR"
$ws.Range("F179").Value = "<value>"

$ws.Range("B180").Value = "SingleUseId349"
$ws.Range("C180").Value = "Large"
$ws.Range("D180").Value = "Left"
$ws.Range("E180").Value = "LTR"
$ws.Range("F180").Value = "None"

$ws.Range("B181").Value = "SingleUseId350"
$ws.Range("C181").Value = "Large"
$ws.Range("D181").Value = "Left"
$ws.Range("E181").Value = "LTR"
$ws.Range("F181").Value = "<value>"

$ws.Range("B182").Value = "SingleUseId351"
$ws.Range("C182").Value = "Large"
$ws.Range("D182").Value = "Left"
$ws.Range("E182").Value = "LTR"
$ws.Range("F182").Value = "None"

$ws.Range("B183").Value = "SingleUseId352"
$ws.Range("C183").Value = "Large"
$ws.Range("D183").Value = "Left"
$ws.Range("E183").Value = "LTR"
$ws.Range("F183").Value = "<value>"

$ws.Range("B184").Value = "SingleUseId353"
$ws.Range("C184").Value = "Large"
$ws.Range("D184").Value = "Left"
$ws.Range("E184").Value = "LTR"
$ws.Range("F184").Value = "None"

$ws.Range("B185").Value = "SingleUseId354"
$ws.Range("C185").Value = "Large"
$ws.Range("D185").Value = "Left"
$ws.Range("E185").Value = "LTR"
$ws.Range("F185").Value = "<value>"

$ws.Range("B186").Value = "SingleUseId355"
$ws.Range("C186").Value = "Large"
$ws.Range("D186").Value = "Left"
$ws.Range("E186").Value = "LTR"
$ws.Range("F186").Value = "None"
